$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be forced to Text
# format first, otherwise Excel auto-converts the string to a numeric value
# (losing the trailing zero / exact text representation used by this sheet).

$ws.Range("D2").Value = "34.007.43"
$ws.Range("E2").Value = "  +10.71%  "

$ws.Range("D3").Value = "1.812.95"
$ws.Range("E3").Value = "  +7.32%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.04"
$ws.Range("E5").Value = "  +3.09%  "

$ws.Range("E6").Value = "  +3.64%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.85"
$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.89"
$ws.Range("E9").Value = "  +3.49%  "

$ws.Range("E10").Value = "  +4.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0668"
$ws.Range("E11").Value = "  +6.69%  "

$ws.Range("E12").Value = "  +2.37%  "

$ws.Range("D13").Value = "2.076.59"
$ws.Range("E13").Value = "  +7.36%  "

$ws.Range("D14").Value = "1.821.31"
$ws.Range("E14").Value = "  +7.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.637"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "33.956.93"
$ws.Range("E16").Value = "  +10.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.12"
$ws.Range("E17").Value = "  -5.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.26"
$ws.Range("E18").Value = "  +6.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.80"
$ws.Range("E19").Value = "  +3.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "254.92"
$ws.Range("E20").Value = "  +3.18%  "

$ws.Range("D21").Value = "0.0₃0740"
$ws.Range("E21").Value = "  +3.38%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.35"
$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.30"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  +1.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.12"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.51"
$ws.Range("E27").Value = "  +3.80%  "

$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.115"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").Value = "  +4.59%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("E31").Value = "  +9.03%  "

$ws.Range("E32").Value = "  +5.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0507"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.50"
$ws.Range("E34").Value = "  +6.25%  "

$ws.Range("D35").Value = "1.550.15"
$ws.Range("E35").Value = "  +2.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.80"
$ws.Range("E36").Value = "  +3.37%  "

$ws.Range("E37").Value = "  +3.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0186"
$ws.Range("E38").Value = "  +3.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.82"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.614"
$ws.Range("E40").Value = "  +4.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.84"
$ws.Range("E41").Value = "  +4.59%  "

$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.899"
$ws.Range("E43").Value = "  +5.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.09"
$ws.Range("E44").Value = "  +4.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0524"
$ws.Range("E45").Value = "  +3.82%  "

$ws.Range("E46").Value = "  +3.07%  "

$ws.Range("D47").Value = "1.963.76"
$ws.Range("E47").Value = "  +7.31%  "

$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.66"
$ws.Range("E49").Value = "  +3.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.73"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "96.13"
$ws.Range("E51").Value = "  +1.10%  "
